# Update the "nota_view" column (J) for the forum grading week 11/09/2022-17/09/2022:
# every row whose nota_view is currently 5 becomes 4 (rows already at 0 stay untouched).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 10).End(-4162).Row  # xlUp

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 10)  # column J = 10
    if ($cell.Value() -eq 5) {
        $cell.Value = 4
    }
}
